$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first two data rows (2008, 2009) - rows 2 and 3.
# This shifts 2010..2020 up to rows 2..12.
$ws.Range("A2:H3").Delete() | Out-Null

# Append the new row for 2021 at row 13 (B column stays blank - no data reported
# for that indicator that year, same pattern as the other rows above).
$ws.Range("A13").Value = "2021年"
$ws.Range("C13").Value = 0.08
$ws.Range("D13").Value = 3506.0478
$ws.Range("E13").Value = 0.14
$ws.Range("F13").Value = 84.72033436
$ws.Range("G13").Value = 80.4135295
$ws.Range("H13").Value = 54873.0348

# Copy the style from the year column (A12 already has style s=1) to the new A13 cell.
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
